# Updates the "Price" (column D) and "Volume(1h)" (column E) cells in the
# cryptos list on the active worksheet, matching the latest scrape.
#
# Numeric-looking Price strings are entered with a leading apostrophe so
# Excel keeps them as literal text (preserving formatting such as trailing
# zeros and multi-dot "thousands.hundreds" groupings) instead of silently
# re-interpreting them as numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.769.90"
$ws.Range("E2").Value = "  +0.12%  "

$ws.Range("D3").Value = "2.082.66"
$ws.Range("E3").Value = "  -1.22%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").Value = "'234.24"
$ws.Range("E5").Value = "  -0.63%  "

$ws.Range("D6").Value = "'0.626"
$ws.Range("E6").Value = "  -0.23%  "

$ws.Range("D7").Value = "'58.63"
$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  -0.03%  "

$ws.Range("D9").Value = "'0.392"
$ws.Range("E9").Value = "  +0.12%  "

$ws.Range("D10").Value = "'0.0787"
$ws.Range("E10").Value = "  +0.35%  "

$ws.Range("E11").Value = "  +2.63%  "

$ws.Range("D12").Value = "2.388.76"
$ws.Range("E12").Value = "  -1.05%  "

$ws.Range("D13").Value = "'14.83"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("D14").Value = "'21.08"
$ws.Range("E14").Value = "  -1.66%  "

$ws.Range("E15").Value = "  -2.40%  "

$ws.Range("D16").Value = "'5.30"
$ws.Range("E16").Value = "  +1.11%  "

$ws.Range("D17").Value = "2.077.14"
$ws.Range("E17").Value = "  -1.16%  "

$ws.Range("D18").Value = "37.689.79"
$ws.Range("E18").Value = "  +0.11%  "

$ws.Range("E19").Value = "  -0.44%  "

$ws.Range("D20").Value = "'71.24"
$ws.Range("E20").Value = "  +1.92%  "

$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("D22").Value = "'228.62"
$ws.Range("E22").Value = "  +0.58%  "

$ws.Range("E23").Value = "  -0.13%  "

$ws.Range("E24").Value = "  -1.43%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").Value = "'169.22"

$ws.Range("E27").Value = "  +3.10%  "

$ws.Range("D28").Value = "'9.00"
$ws.Range("E28").Value = "  +0.39%  "

$ws.Range("D29").Value = "'19.52"
$ws.Range("E29").Value = "  +0.14%  "

$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("D31").Value = "'0.121"
$ws.Range("E31").Value = "  +1.57%  "

$ws.Range("E32").Value = "  +0.54%  "

$ws.Range("D33").Value = "'0.0632"
$ws.Range("E33").Value = "  +1.18%  "

$ws.Range("E34").Value = "  +1.23%  "

$ws.Range("E35").Value = "  -3.46%  "

$ws.Range("E36").Value = "  +2.25%  "

$ws.Range("D37").Value = "'3.39"
$ws.Range("E37").Value = "  -4.64%  "

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = "  -0.04%  "

$ws.Range("D39").Value = "'5.39"
$ws.Range("E39").Value = "  -4.33%  "

$ws.Range("D40").Value = "'0.0978"
$ws.Range("E40").Value = "  +1.38%  "

$ws.Range("D41").Value = "'98.32"
$ws.Range("E41").Value = "  +0.31%  "

$ws.Range("E42").Value = "  +0.68%  "

$ws.Range("E43").Value = "  -2.95%  "

$ws.Range("D44").Value = "1.459.55"
$ws.Range("E44").Value = "  -1.85%  "

$ws.Range("E45").Value = "  +4.38%  "

$ws.Range("E46").Value = "  -0.38%  "

$ws.Range("D47").Value = "'16.49"
$ws.Range("E47").Value = "  +5.85%  "

$ws.Range("E48").Value = "  +0.87%  "

$ws.Range("D49").Value = "'7.43"
$ws.Range("E49").Value = "  +1.35%  "

$ws.Range("D50").Value = "'3.02"
$ws.Range("E50").Value = "  -0.36%  "

$ws.Range("D51").Value = "2.273.76"
$ws.Range("E51").Value = "  -1.22%  "
